# Revert file upload functionality
# Re-add row 59 (previously removed) to each of the 4 data sheets.

$wb = $excel.ActiveWorkbook

$rows = @{
    1 = @{
        A = 45845.46043981481
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x64"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 356
        I = 7
    }
    2 = @{
        A = 45845.46043981481
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x5C"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 348
        I = 25
    }
    3 = @{
        A = 45845.46043981481
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x68"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 104
        I = 15
    }
    4 = @{
        A = 45845.46043981481
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7D"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 125
        I = 9
    }
}

foreach ($idx in 1,2,3,4) {
    $ws = $wb.Worksheets.Item($idx)
    $data = $rows[$idx]

    $ws.Range("A59").Value = $data.A
    $ws.Range("B59").Value = $data.B
    $ws.Range("C59").Value = $data.C
    $ws.Range("D59").Value = $data.D
    $ws.Range("E59").Value = $data.E
    $ws.Range("F59").Value = $data.F
    $ws.Range("G59").Value = [double]$data.G
    $ws.Range("H59").Value = $data.H
    $ws.Range("I59").Value = $data.I

    # Match the date/time number format used by the cell above it (A58).
    $ws.Range("A59").NumberFormat = $ws.Range("A58").NumberFormat()
}
